$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Capture the plain "list" style (currently on A2) before we overwrite A2,
#    so we can re-apply it further down the column (A13:A15).
$ws.Range("A2").Copy()
$ws.Range("A13:A15").PasteSpecial(-4122)

# 2) Write the new node ids into A2:A5 (replacing the old comma-joined text
#    that lived in A2, and filling the previously-blank A3:A5).
$ws.Range("A2").Value = 146568
$ws.Range("A3").Value = 202047
$ws.Range("A4").Value = 215069
$ws.Range("A5").Value = 279290

# 3) A2:A4 pick up the "NODES" block formatting (same as A3 originally had),
#    then get a lighter top border to separate them from the header row.
$ws.Range("A3").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)
$ws.Range("A2:A4").Borders.Item(7).Color = 0

# 4) A5:A12 pick up the alternate "NODES" block formatting already used by B8.
$ws.Range("B8").Copy()
$ws.Range("A5:A12").PasteSpecial(-4122)

# 5) C2:C3 pick up the plain block formatting used throughout column A/B.
$ws.Range("A16").Copy()
$ws.Range("C2:C3").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# 6) Selection ends on A6.
$ws.Range("A6").Select()
